$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,9,12,14,16,17,18,21,23,24,25,26,29,31,33,35,36,37,38,39,42,48,52,56,58,59,61,63,64,66,76,77,78,79,81,82,84,85,86,88,90,92,93,95,97,99,100,101,104,105,106,107,109,113,114,122,128,130,132,133,138,141,143,148,149,154,156,166,167,168,170,175,182,184,185,188,190,193,198,199,202,206,209,211,224,226,227,229,232)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 15).Value = ""
}
